# Update TPM-derived values in the LR-pair sheet (Edn3-Ednrb)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1780495
$ws.Range("H2").Value = 0.356099
$ws.Range("M2").Value = 90.353905
$ws.Range("N2").Value = 180.70781
$ws.Range("O2").Value = 0.3131638580342592
$ws.Range("P2").Value = 0.2965570309229201
$ws.Range("Q2").Value = 16.0874676082975
$ws.Range("R2").Value = 64.34987043318999
$ws.Range("S2").Value = 0.3131638580342592
$ws.Range("T2").Value = 0.2965570309229201

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1780495
$ws.Range("H3").Value = 0.356099
$ws.Range("O3").Value = 0.001245676287098259
$ws.Range("P3").Value = 0.001769428615638018
$ws.Range("Q3").Value = 0.063991346399
$ws.Range("R3").Value = 0.383948078394
$ws.Range("S3").Value = 0.001245676287098259
$ws.Range("T3").Value = 0.001769428615638018

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.1780495
$ws.Range("H4").Value = 0.356099
$ws.Range("M4").Value = 10.80810533333333
$ws.Range("N4").Value = 32.424316
$ws.Range("O4").Value = 0.03746056093787335
$ws.Range("P4").Value = 0.05321108635352579
$ws.Range("Q4").Value = 1.924377750547333
$ws.Range("R4").Value = 11.546266503284
$ws.Range("S4").Value = 0.03746056093787335
$ws.Range("T4").Value = 0.05321108635352579

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.1780495
$ws.Range("H5").Value = 0.356099
$ws.Range("M5").Value = 165.852196
$ws.Range("N5").Value = 331.704392
$ws.Range("O5").Value = 0.574838614477306
$ws.Range("P5").Value = 0.5443553858331436
$ws.Range("Q5").Value = 29.529900571702
$ws.Range("R5").Value = 118.119602286808
$ws.Range("S5").Value = 0.574838614477306
$ws.Range("T5").Value = 0.5443553858331436

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.1780495
$ws.Range("H6").Value = 0.356099
$ws.Range("M6").Value = 5.448456
$ws.Range("N6").Value = 16.345368
$ws.Range("O6").Value = 0.01888418105769649
$ws.Range("P6").Value = 0.02682415222353981
$ws.Range("Q6").Value = 0.9700948665720001
$ws.Range("R6").Value = 5.820569199432001
$ws.Range("S6").Value = 0.01888418105769649
$ws.Range("T6").Value = 0.02682415222353981

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.1780495
$ws.Range("H7").Value = 0.356099
$ws.Range("M7").Value = 15.69751633333333
$ws.Range("N7").Value = 47.09254900000001
$ws.Range("O7").Value = 0.05440710920576665
$ws.Range("P7").Value = 0.07728291605123282
$ws.Range("Q7").Value = 2.794934934391834
$ws.Range("R7").Value = 16.769609606351
$ws.Range("S7").Value = 0.05440710920576665
$ws.Range("T7").Value = 0.07728291605123282
